$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 11: E11 15->16, F11 10->11, H11 11->12
$ws.Range("E11").Value = 16
$ws.Range("F11").Value = 11
$ws.Range("H11").Value = 12

# Row 12: E12 27->28
$ws.Range("E12").Value = 28
